# Apply cryptos-list price/volume refresh (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "87.346.08"
$ws.Range("E2").Value = "  -1.90%  "
$ws.Range("D3").Value = "3.160.72"
$ws.Range("E3").Value = "  -6.66%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'204.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.57%  "
$ws.Range("D6").Value = "'608.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.68%  "
$ws.Range("E7").Value = "  -8.21%  "
$ws.Range("D8").Value = "'0.661"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.92%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "3.151.06"
$ws.Range("E10").Value = "  -6.78%  "
$ws.Range("D11").Value = "'0.533"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -16.36%  "
$ws.Range("D12").Value = "'0.178"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.28%  "
$ws.Range("E13").Value = "  -16.49%  "
$ws.Range("D14").Value = "3.743.64"
$ws.Range("E14").Value = "  -6.32%  "
$ws.Range("E15").Value = "  -5.89%  "
$ws.Range("D16").Value = "87.129.38"
$ws.Range("E16").Value = "  -2.08%  "
$ws.Range("D17").Value = "'31.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -13.02%  "
$ws.Range("D18").Value = "3.185.23"
$ws.Range("E18").Value = "  -4.78%  "
$ws.Range("D19").Value = "'2.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.12%  "
$ws.Range("E20").Value = "  -10.60%  "
$ws.Range("D21").Value = "'413.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.93%  "
$ws.Range("D22").Value = "'8.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -13.13%  "
$ws.Range("D23").Value = "'5.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -10.01%  "
$ws.Range("D24").Value = "'5.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.96%  "
$ws.Range("D25").Value = "'11.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.64%  "
$ws.Range("D26").Value = "3.323.86"
$ws.Range("E26").Value = "  -6.36%  "
$ws.Range("D27").Value = "'73.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.06%  "
$ws.Range("D28").Value = "'0.0000129"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -10.12%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "'0.160"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -19.51%  "
$ws.Range("D31").Value = "'0.997"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").Value = "'537.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -10.18%  "
$ws.Range("D33").Value = "'8.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -13.95%  "
$ws.Range("D34").Value = "'1.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -17.91%  "
$ws.Range("E35").Value = "  -9.15%  "
$ws.Range("D36").Value = "'1.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -13.15%  "
$ws.Range("E37").Value = "  -8.68%  "
$ws.Range("D38").Value = "'21.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.87%  "
$ws.Range("D39").Value = "'21.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("D40").Value = "'0.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").Value = "'2.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.21%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("E43").Value = "  -13.23%  "
$ws.Range("D44").Value = "'0.370"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -14.61%  "
$ws.Range("D45").Value = "'147.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.19%  "
$ws.Range("D46").Value = "'172.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.13%  "
$ws.Range("D47").Value = "'43.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.97%  "
$ws.Range("D48").Value = "'0.126"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.57%  "
$ws.Range("D49").Value = "'1.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -15.66%  "
$ws.Range("D50").Value = "'3.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -12.71%  "
$ws.Range("D51").Value = "'0.694"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -12.23%  "
